# Updates the cryptos price/volume table (Coin, Link, Price, Volume(1h)) to the
# latest scrape snapshot. Mirrors the GitHub Actions refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.964.57'
$ws.Range('E2').Value = '  -0.58%  '
# Row 3
$ws.Range('D3').Value = '1.826.21'
# Row 4
$ws.Range('E4').Value = '  -0.32%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.29'
$ws.Range('E5').Value = '  +0.07%  '
# Row 6
$ws.Range('E6').Value = '  -0.25%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4573'
$ws.Range('E7').Value = '  -1.23%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3699'
$ws.Range('E8').Value = '  +1.87%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07307'
$ws.Range('E9').Value = '  +0.03%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8722'
$ws.Range('E10').Value = '  +0.18%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07960'
$ws.Range('E11').Value = '  +4.54%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.74'
$ws.Range('E12').Value = '  -1.95%  '
# Row 13
$ws.Range('D13').Value = '1.854.67'
$ws.Range('E13').Value = '  +1.03%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.568'
$ws.Range('E14').Value = '  +1.45%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.321'
$ws.Range('E15').Value = '  -0.36%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.33'
$ws.Range('E16').Value = '  -1.19%  '
# Row 17
$ws.Range('E17').Value = '  -0.11%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008848'
$ws.Range('E18').Value = '  +2.23%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.009'
$ws.Range('E19').Value = '  -0.05%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.69'
$ws.Range('E20').Value = '  +1.45%  '
# Row 21
$ws.Range('D21').Value = '27.329.96'
$ws.Range('E21').Value = '  -0.51%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.094'
$ws.Range('E22').Value = '  -2.38%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.53'
$ws.Range('E23').Value = '  -0.25%  '
# Row 24
$ws.Range('D24').Value = '2.146.99'
$ws.Range('E24').Value = '  +2.18%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.18'
$ws.Range('E25').Value = '  +1.18%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.847'
$ws.Range('E26').Value = '  -0.88%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.37'
$ws.Range('E27').Value = '  +1.15%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.039'
$ws.Range('E28').Value = '  -1.81%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.140'
$ws.Range('E29').Value = '  +0.84%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '114.96'
$ws.Range('E30').Value = '  -0.93%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08874'
$ws.Range('E31').Value = '  -0.33%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.957'
$ws.Range('E32').Value = '  -0.08%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7277'
$ws.Range('E33').Value = '  -1.35%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.441'
$ws.Range('E34').Value = '  -0.28%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.127'
$ws.Range('E35').Value = '  -1.06%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.072'
$ws.Range('E36').Value = '  +0.21%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.435'
$ws.Range('E37').Value = '  -1.71%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01939'
$ws.Range('E38').Value = '  +1.20%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05220'
$ws.Range('E39').Value = '  -0.58%  '
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.935'
$ws.Range('E40').Value = '  +0.39%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.152'
$ws.Range('E41').Value = '  -0.33%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5131'
$ws.Range('E42').Value = '  -1.41%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1627'
$ws.Range('E43').Value = '  -0.13%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.174'
$ws.Range('E44').Value = '  -1.32%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4820'
$ws.Range('E45').Value = '  -0.48%  '
# Row 46
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.007'
$ws.Range('E46').Value = '  -0.25%  '
# Row 47
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.16'
$ws.Range('E47').Value = '  +0.21%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.42'
$ws.Range('E48').Value = '  -0.88%  '
# Row 49
$ws.Range('E49').Value = '  -0.55%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06209'
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '64.90'
$ws.Range('E51').Value = '  +0.63%  '
